# Apply the "new account statement" edit:
#  - Insert a new worker data row (row 17) below the existing one (row 16)
#    with the new worker's data, formatted like the existing data row.
#  - Update the aggregate totals (VALOR MORA, Cant. Trabajadores, Cant. Periodos)
#  - The signature block rows shift down by one row to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 17 (pushes old row17.. down by one, fixes merged
#    cells and the rows below automatically).
$ws.Rows("17").Insert()

# 2. Copy the formatting of the existing data row (16) onto the new row (17)
#    so it keeps the same borders/fonts/number formats.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3. Fill in the new worker's data on row 17.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1051885632"
$ws.Range("D17").Value = "ANA MILENA VASQUEZ BELTRAN"
$ws.Range("E17").Value = "1801"
$ws.Range("F17").Value = 3935
$ws.Range("G17").Value = 737717

# 4. Update the aggregate summary figures.
$ws.Range("E11").Value = 15739
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2
